$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = 45830
$ws.Range("B82").Value = 30.26
$ws.Range("C82").Value = 89.5
$ws.Range("D82").Formula = "=B82-B81"
$ws.Range("E82").Formula = "=C82-C81"

$ws.Range("E83").Select()
